# Apply playtesting-feedback updates to the Sprint Plan workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing task status tweaks -------------------------------------------------
# Row 7 "Create game pitch presentation" -> Finished? WIP -> Yes
$ws.Range("E7").Value = "Yes"

# Row 27 "Begin playtesting" -> Finished? (blank) -> Yes
$ws.Range("E27").Value = "Yes"

# Row 39 "Make design changes based on playtesting feedback" -> Finished? (blank) -> WIP
$ws.Range("E39").Value = "WIP"

# --- New "PLAYTESTING FEEDBACK GOALS" list in columns G:I -----------------------
$ws.Range("G16").Value = "PLAYTESTING FEEDBACK GOALS"

$ws.Range("G17").Value = "Add knock back to hero"
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = "Yes"

$ws.Range("G18").Value = "Add pathfinding to ghosts"
$ws.Range("H18").Value = 1

$ws.Range("G19").Value = "Make camera speed proportional to hero speed"
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = "Yes"

$ws.Range("G20").Value = "Make sanity/hp reset per level"
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = "Yes"

$ws.Range("G21").Value = "Make sanity wobble more punishing"
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = "Yes"

$ws.Range("G22").Value = "Make UI more obvious"
$ws.Range("H22").Value = 1

$ws.Range("G23").Value = "Invest in shaders"
$ws.Range("H23").Value = 3

$ws.Range("G24").Value = "Noise indicator for taking damage"
$ws.Range("H24").Value = 1

$ws.Range("G25").Value = "Pause spawners while at altars"
$ws.Range("H25").Value = 1
$ws.Range("I25").Value = "Yes"

$ws.Range("G26").Value = "Fix font in sanity messages"
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = "Yes"

# --- View state (scrolled/selected cell) matches the author's last look --------
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("I24").Select()
